$wb = $excel.ActiveWorkbook

# ===================== Common 'F' (want-to-go count) tweaks, rows 2-22 =====================
# Identical edits on sheet 1 (展览) and sheet 4 (全部类型).
foreach ($idx in @(1,4)) {
    $s = $wb.Worksheets.Item($idx)
    $s.Cells.Item(2, 6).Value = 3160
    $s.Cells.Item(3, 6).Value = 562
    $s.Cells.Item(5, 6).Value = 120
    $s.Cells.Item(6, 6).Value = 72
    $s.Cells.Item(9, 6).Value = 1161
    $s.Cells.Item(10, 6).Value = 16353
    $s.Cells.Item(11, 6).Value = 277
    $s.Cells.Item(12, 6).Value = 203
    $s.Cells.Item(14, 6).Value = 6364
    $s.Cells.Item(15, 6).Value = 638
    $s.Cells.Item(17, 6).Value = 81
    $s.Cells.Item(18, 6).Value = 22
    $s.Cells.Item(20, 6).Value = 1272
    $s.Cells.Item(21, 6).Value = 41
    $s.Cells.Item(22, 6).Value = 37
}

# ===================== Sheet 1: 展览 =====================
# ---- ws1: Worksheets.Item(1) ----
$ws1 = $wb.Worksheets.Item(1)

# A brand-new row is inserted at row 26 (a newly announced event); every
# row that used to sit at 26..old_last_row now holds the data of the row
# immediately above it, and the sheet gains one new row at the bottom that
# repeats the old final row's data. Column A (the plain row index) is left
# as-is for the rows that only shift content, and is newly written only for
# the freshly created last row.

# give the brand-new final row (A40) the same style as A39 (bold/boxed index cell)
$ws1.Cells.Item(39, 1).Copy($ws1.Cells.Item(40, 1))
$ws1.Cells.Item(40, 1).Value = 39

# row 26: 苏州·第一届维度创想动漫嘉年华
$ws1.Cells.Item(26, 2).NumberFormat = "@"
$ws1.Cells.Item(26, 2).Value = '2024-07-28'
$ws1.Cells.Item(26, 2).ClearFormats()
$ws1.Cells.Item(26, 3).Value = '苏州·第一届维度创想动漫嘉年华'
$ws1.Cells.Item(26, 4).Value = '广济北路799号 平江市民健身中心'
$ws1.Cells.Item(26, 5).Value = '2024.07.28 09:00-07.28 18:00'
$ws1.Cells.Item(26, 6).Value = 0
$ws1.Cells.Item(26, 7).Value = 60
$ws1.Cells.Item(26, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88083'
$ws1.Cells.Item(26, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/j2YxrXMV1719311394987.jpeg'

# row 27: 苏州·环球港动漫节
$ws1.Cells.Item(27, 2).NumberFormat = "@"
$ws1.Cells.Item(27, 2).Value = '2024-08-02'
$ws1.Cells.Item(27, 2).ClearFormats()
$ws1.Cells.Item(27, 3).Value = '苏州·环球港动漫节'
$ws1.Cells.Item(27, 4).Value = '相城大道1609号 苏州环球港'
$ws1.Cells.Item(27, 5).Value = '2024.08.02 10:00-08.04 16:00'
$ws1.Cells.Item(27, 6).Value = 15
$ws1.Cells.Item(27, 7).Value = 49
$ws1.Cells.Item(27, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87768'
$ws1.Cells.Item(27, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/U6aDJ0Yq1718876195019.jpeg'

# row 28: 常熟·ACG动漫游戏嘉年华
$ws1.Cells.Item(28, 2).NumberFormat = "@"
$ws1.Cells.Item(28, 2).Value = '2024-08-03'
$ws1.Cells.Item(28, 2).ClearFormats()
$ws1.Cells.Item(28, 3).Value = '常熟·ACG动漫游戏嘉年华'
$ws1.Cells.Item(28, 4).Value = '冬青路88号 江南·美好汇生活广场'
$ws1.Cells.Item(28, 5).Value = '2024.08.03 09:00-08.04 17:00'
$ws1.Cells.Item(28, 6).Value = 221
$ws1.Cells.Item(28, 7).Value = 60
$ws1.Cells.Item(28, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85851'
$ws1.Cells.Item(28, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/LgJRjcDn1715933608635.jpeg'

# row 29: 常熟·CDW.动漫展05
$ws1.Cells.Item(29, 2).NumberFormat = "@"
$ws1.Cells.Item(29, 2).Value = '2024-08-03'
$ws1.Cells.Item(29, 2).ClearFormats()
$ws1.Cells.Item(29, 3).Value = '常熟·CDW.动漫展05'
$ws1.Cells.Item(29, 4).Value = '开元大道1号 常熟国际博览中心'
$ws1.Cells.Item(29, 5).Value = '2024.08.03 09:00-08.04 17:00'
$ws1.Cells.Item(29, 6).Value = 895
$ws1.Cells.Item(29, 7).Value = 60
$ws1.Cells.Item(29, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86329'
$ws1.Cells.Item(29, 9).Value = '//i0.hdslb.com/bfs/openplatform/202405/GluYLcFY1716136984305.jpeg'

# row 30: 苏州·代号鸢only茶话会-星渡咖啡
$ws1.Cells.Item(30, 2).NumberFormat = "@"
$ws1.Cells.Item(30, 2).Value = '2024-08-03'
$ws1.Cells.Item(30, 2).ClearFormats()
$ws1.Cells.Item(30, 3).Value = '苏州·代号鸢only茶话会-星渡咖啡'
$ws1.Cells.Item(30, 4).Value = '德必姑苏WE国际文化艺术中心6-102室渔郎桥浜路16号 星渡咖啡'
$ws1.Cells.Item(30, 5).Value = '2024.08.03 10:00-08.04 19:00'
$ws1.Cells.Item(30, 6).Value = 56
$ws1.Cells.Item(30, 7).Value = 50
$ws1.Cells.Item(30, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87685'
$ws1.Cells.Item(30, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/eyHRVQuv1718780132754.jpeg'

# row 31: 苏州·星部落动漫嘉年华
$ws1.Cells.Item(31, 2).NumberFormat = "@"
$ws1.Cells.Item(31, 2).Value = '2024-08-03'
$ws1.Cells.Item(31, 2).ClearFormats()
$ws1.Cells.Item(31, 3).Value = '苏州·星部落动漫嘉年华'
$ws1.Cells.Item(31, 4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Cells.Item(31, 5).Value = '2024.08.03 09:00-08.04 16:00'
$ws1.Cells.Item(31, 6).Value = 5051
$ws1.Cells.Item(31, 7).Value = 68
$ws1.Cells.Item(31, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84858'
$ws1.Cells.Item(31, 9).Value = '//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg'

# row 32: 苏州·授渔仲夏动漫节2.0
$ws1.Cells.Item(32, 2).NumberFormat = "@"
$ws1.Cells.Item(32, 2).Value = '2024-08-04'
$ws1.Cells.Item(32, 2).ClearFormats()
$ws1.Cells.Item(32, 3).Value = '苏州·授渔仲夏动漫节2.0'
$ws1.Cells.Item(32, 4).Value = '一干河东路333号 张家港沙洲湖酒店'
$ws1.Cells.Item(32, 5).Value = '2024.08.04 09:30-08.04 16:30'
$ws1.Cells.Item(32, 6).Value = 500
$ws1.Cells.Item(32, 7).Value = 40
$ws1.Cells.Item(32, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87156'
$ws1.Cells.Item(32, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/31htgxuC1718083875401.jpeg'

# row 33: 苏州·ICAN summer World动漫品牌夏游节
$ws1.Cells.Item(33, 2).NumberFormat = "@"
$ws1.Cells.Item(33, 2).Value = '2024-08-17'
$ws1.Cells.Item(33, 2).ClearFormats()
$ws1.Cells.Item(33, 3).Value = '苏州·ICAN summer World动漫品牌夏游节'
$ws1.Cells.Item(33, 4).Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws1.Cells.Item(33, 5).Value = '2024.08.17 10:00-08.18 17:00'
$ws1.Cells.Item(33, 6).Value = 11331
$ws1.Cells.Item(33, 7).Value = 60
$ws1.Cells.Item(33, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85289'
$ws1.Cells.Item(33, 9).Value = '//i0.hdslb.com/bfs/openplatform/202404/JavlW9fj1714459472747.jpeg'

# row 34: 苏州·第二届Redamancy动漫游戏嘉年华
$ws1.Cells.Item(34, 2).NumberFormat = "@"
$ws1.Cells.Item(34, 2).Value = '2024-08-17'
$ws1.Cells.Item(34, 2).ClearFormats()
$ws1.Cells.Item(34, 3).Value = '苏州·第二届Redamancy动漫游戏嘉年华'
$ws1.Cells.Item(34, 4).Value = '清禾路886号 尹山湖大剧院'
$ws1.Cells.Item(34, 5).Value = '2024.08.17 10:00-08.18 17:00'
$ws1.Cells.Item(34, 6).Value = 1247
$ws1.Cells.Item(34, 7).Value = 60
$ws1.Cells.Item(34, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83576'
$ws1.Cells.Item(34, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/UbwYg1jn1716516632746.jpeg'

# row 35: 苏州·赛马娘ONLY
$ws1.Cells.Item(35, 2).NumberFormat = "@"
$ws1.Cells.Item(35, 2).Value = '2024-08-24'
$ws1.Cells.Item(35, 2).ClearFormats()
$ws1.Cells.Item(35, 3).Value = '苏州·赛马娘ONLY'
$ws1.Cells.Item(35, 4).Value = '东苑路115-11号 苏苑饭店'
$ws1.Cells.Item(35, 5).Value = '2024.08.24 10:00-08.24 16:00'
$ws1.Cells.Item(35, 6).Value = 18
$ws1.Cells.Item(35, 7).Value = 60
$ws1.Cells.Item(35, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87619'
$ws1.Cells.Item(35, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/2seg6U5l1718215220516.jpeg'

# row 36: 苏州·Good jump ACG中秋嘉年华动漫国潮文化节
$ws1.Cells.Item(36, 2).NumberFormat = "@"
$ws1.Cells.Item(36, 2).Value = '2024-09-15'
$ws1.Cells.Item(36, 2).ClearFormats()
$ws1.Cells.Item(36, 3).Value = '苏州·Good jump ACG中秋嘉年华动漫国潮文化节'
$ws1.Cells.Item(36, 4).Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws1.Cells.Item(36, 5).Value = '2024.09.15 10:00-09.16 17:00'
$ws1.Cells.Item(36, 6).Value = 150
$ws1.Cells.Item(36, 7).Value = 60
$ws1.Cells.Item(36, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87120'
$ws1.Cells.Item(36, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/5Qp6CuZ11717828215727.jpeg'

# row 37: 苏州·I COME ACG动漫品牌博览会
$ws1.Cells.Item(37, 2).NumberFormat = "@"
$ws1.Cells.Item(37, 2).Value = '2024-10-01'
$ws1.Cells.Item(37, 2).ClearFormats()
$ws1.Cells.Item(37, 3).Value = '苏州·I COME ACG动漫品牌博览会'
$ws1.Cells.Item(37, 4).Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws1.Cells.Item(37, 5).Value = '2024.10.01 10:00-10.03 17:00'
$ws1.Cells.Item(37, 6).Value = 207
$ws1.Cells.Item(37, 7).Value = 60
$ws1.Cells.Item(37, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87118'
$ws1.Cells.Item(37, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/SyK3KnNb1717829071115.jpeg'

# row 38: 苏州·第十三届理想乡动漫展-同人创作者大会
$ws1.Cells.Item(38, 2).NumberFormat = "@"
$ws1.Cells.Item(38, 2).Value = '2024-10-01'
$ws1.Cells.Item(38, 2).ClearFormats()
$ws1.Cells.Item(38, 3).Value = '苏州·第十三届理想乡动漫展-同人创作者大会'
$ws1.Cells.Item(38, 4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws1.Cells.Item(38, 5).Value = '2024.10.01 10:00-10.03 17:00'
$ws1.Cells.Item(38, 6).Value = 3838
$ws1.Cells.Item(38, 7).Value = 39
$ws1.Cells.Item(38, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83821'
$ws1.Cells.Item(38, 9).Value = '//i0.hdslb.com/bfs/openplatform/202404/OMtuTTFY1711958198579.jpeg'

# row 39: 苏州·明日方舟ONLY#2024~佑桑柔
$ws1.Cells.Item(39, 2).NumberFormat = "@"
$ws1.Cells.Item(39, 2).Value = '2024-10-02'
$ws1.Cells.Item(39, 2).ClearFormats()
$ws1.Cells.Item(39, 3).Value = '苏州·明日方舟ONLY#2024~佑桑柔'
$ws1.Cells.Item(39, 4).Value = '城际路21号 苏州汇融广场假日酒店'
$ws1.Cells.Item(39, 5).Value = '2024.10.02 10:00-10.02 17:00'
$ws1.Cells.Item(39, 6).Value = 270
$ws1.Cells.Item(39, 7).Value = 75
$ws1.Cells.Item(39, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84046'
$ws1.Cells.Item(39, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/0VhIRprD1716344515303.jpeg'

# row 40: 苏州·第三届华盟国漫次元嘉年华
$ws1.Cells.Item(40, 2).NumberFormat = "@"
$ws1.Cells.Item(40, 2).Value = '2024-10-26'
$ws1.Cells.Item(40, 2).ClearFormats()
$ws1.Cells.Item(40, 3).Value = '苏州·第三届华盟国漫次元嘉年华'
$ws1.Cells.Item(40, 4).Value = '清禾路886号 苏州聚橙尹山湖大剧院'
$ws1.Cells.Item(40, 5).Value = '2024.10.26 10:00-10.27 17:00'
$ws1.Cells.Item(40, 6).Value = 74
$ws1.Cells.Item(40, 7).Value = 58
$ws1.Cells.Item(40, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85767'
$ws1.Cells.Item(40, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/CqSYBZhQ1715846719965.jpeg'

# ===================== Sheet 4: 全部类型 =====================
# ---- ws4: Worksheets.Item(4) ----
$ws4 = $wb.Worksheets.Item(4)

# A brand-new row is inserted at row 26 (a newly announced event); every
# row that used to sit at 26..old_last_row now holds the data of the row
# immediately above it, and the sheet gains one new row at the bottom that
# repeats the old final row's data. Column A (the plain row index) is left
# as-is for the rows that only shift content, and is newly written only for
# the freshly created last row.

# give the brand-new final row (A41) the same style as A40 (bold/boxed index cell)
$ws4.Cells.Item(40, 1).Copy($ws4.Cells.Item(41, 1))
$ws4.Cells.Item(41, 1).Value = 40

# row 26: 苏州·第一届维度创想动漫嘉年华
$ws4.Cells.Item(26, 2).NumberFormat = "@"
$ws4.Cells.Item(26, 2).Value = '2024-07-28'
$ws4.Cells.Item(26, 2).ClearFormats()
$ws4.Cells.Item(26, 3).Value = '苏州·第一届维度创想动漫嘉年华'
$ws4.Cells.Item(26, 4).Value = '广济北路799号 平江市民健身中心'
$ws4.Cells.Item(26, 5).Value = '2024.07.28 09:00-07.28 18:00'
$ws4.Cells.Item(26, 6).Value = 0
$ws4.Cells.Item(26, 7).Value = 60
$ws4.Cells.Item(26, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=88083'
$ws4.Cells.Item(26, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/j2YxrXMV1719311394987.jpeg'

# row 27: 苏州·环球港动漫节
$ws4.Cells.Item(27, 2).NumberFormat = "@"
$ws4.Cells.Item(27, 2).Value = '2024-08-02'
$ws4.Cells.Item(27, 2).ClearFormats()
$ws4.Cells.Item(27, 3).Value = '苏州·环球港动漫节'
$ws4.Cells.Item(27, 4).Value = '相城大道1609号 苏州环球港'
$ws4.Cells.Item(27, 5).Value = '2024.08.02 10:00-08.04 16:00'
$ws4.Cells.Item(27, 6).Value = 15
$ws4.Cells.Item(27, 7).Value = 49
$ws4.Cells.Item(27, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87768'
$ws4.Cells.Item(27, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/U6aDJ0Yq1718876195019.jpeg'

# row 28: 常熟·ACG动漫游戏嘉年华
$ws4.Cells.Item(28, 2).NumberFormat = "@"
$ws4.Cells.Item(28, 2).Value = '2024-08-03'
$ws4.Cells.Item(28, 2).ClearFormats()
$ws4.Cells.Item(28, 3).Value = '常熟·ACG动漫游戏嘉年华'
$ws4.Cells.Item(28, 4).Value = '冬青路88号 江南·美好汇生活广场'
$ws4.Cells.Item(28, 5).Value = '2024.08.03 09:00-08.04 17:00'
$ws4.Cells.Item(28, 6).Value = 221
$ws4.Cells.Item(28, 7).Value = 60
$ws4.Cells.Item(28, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85851'
$ws4.Cells.Item(28, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/LgJRjcDn1715933608635.jpeg'

# row 29: 常熟·CDW.动漫展05
$ws4.Cells.Item(29, 2).NumberFormat = "@"
$ws4.Cells.Item(29, 2).Value = '2024-08-03'
$ws4.Cells.Item(29, 2).ClearFormats()
$ws4.Cells.Item(29, 3).Value = '常熟·CDW.动漫展05'
$ws4.Cells.Item(29, 4).Value = '开元大道1号 常熟国际博览中心'
$ws4.Cells.Item(29, 5).Value = '2024.08.03 09:00-08.04 17:00'
$ws4.Cells.Item(29, 6).Value = 895
$ws4.Cells.Item(29, 7).Value = 60
$ws4.Cells.Item(29, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86329'
$ws4.Cells.Item(29, 9).Value = '//i0.hdslb.com/bfs/openplatform/202405/GluYLcFY1716136984305.jpeg'

# row 30: 苏州·代号鸢only茶话会-星渡咖啡
$ws4.Cells.Item(30, 2).NumberFormat = "@"
$ws4.Cells.Item(30, 2).Value = '2024-08-03'
$ws4.Cells.Item(30, 2).ClearFormats()
$ws4.Cells.Item(30, 3).Value = '苏州·代号鸢only茶话会-星渡咖啡'
$ws4.Cells.Item(30, 4).Value = '德必姑苏WE国际文化艺术中心6-102室渔郎桥浜路16号 星渡咖啡'
$ws4.Cells.Item(30, 5).Value = '2024.08.03 10:00-08.04 19:00'
$ws4.Cells.Item(30, 6).Value = 56
$ws4.Cells.Item(30, 7).Value = 50
$ws4.Cells.Item(30, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87685'
$ws4.Cells.Item(30, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/eyHRVQuv1718780132754.jpeg'

# row 31: 苏州·星部落动漫嘉年华
$ws4.Cells.Item(31, 2).NumberFormat = "@"
$ws4.Cells.Item(31, 2).Value = '2024-08-03'
$ws4.Cells.Item(31, 2).ClearFormats()
$ws4.Cells.Item(31, 3).Value = '苏州·星部落动漫嘉年华'
$ws4.Cells.Item(31, 4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Cells.Item(31, 5).Value = '2024.08.03 09:00-08.04 16:00'
$ws4.Cells.Item(31, 6).Value = 5051
$ws4.Cells.Item(31, 7).Value = 68
$ws4.Cells.Item(31, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84858'
$ws4.Cells.Item(31, 9).Value = '//i0.hdslb.com/bfs/openplatform/202404/UI5EFZTT1713685680462.jpeg'

# row 32: 苏州·授渔仲夏动漫节2.0
$ws4.Cells.Item(32, 2).NumberFormat = "@"
$ws4.Cells.Item(32, 2).Value = '2024-08-04'
$ws4.Cells.Item(32, 2).ClearFormats()
$ws4.Cells.Item(32, 3).Value = '苏州·授渔仲夏动漫节2.0'
$ws4.Cells.Item(32, 4).Value = '一干河东路333号 张家港沙洲湖酒店'
$ws4.Cells.Item(32, 5).Value = '2024.08.04 09:30-08.04 16:30'
$ws4.Cells.Item(32, 6).Value = 500
$ws4.Cells.Item(32, 7).Value = 40
$ws4.Cells.Item(32, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87156'
$ws4.Cells.Item(32, 9).Value = '//i1.hdslb.com/bfs/openplatform/202406/31htgxuC1718083875401.jpeg'

# row 33: 苏州·爱乐之城·经典电影作品音乐会
$ws4.Cells.Item(33, 2).NumberFormat = "@"
$ws4.Cells.Item(33, 2).Value = '2024-08-10'
$ws4.Cells.Item(33, 2).ClearFormats()
$ws4.Cells.Item(33, 3).Value = '苏州·爱乐之城·经典电影作品音乐会'
$ws4.Cells.Item(33, 4).Value = '念珠街121号道前街与吉庆街路口距养育巷地铁站 苏州市会议中心'
$ws4.Cells.Item(33, 5).Value = '2024.08.10 19:30-08.10 21:00'
$ws4.Cells.Item(33, 6).Value = 21
$ws4.Cells.Item(33, 7).Value = 50
$ws4.Cells.Item(33, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=86194'
$ws4.Cells.Item(33, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/vagzbfox1716438290025.jpeg'

# row 34: 苏州·ICAN summer World动漫品牌夏游节
$ws4.Cells.Item(34, 2).NumberFormat = "@"
$ws4.Cells.Item(34, 2).Value = '2024-08-17'
$ws4.Cells.Item(34, 2).ClearFormats()
$ws4.Cells.Item(34, 3).Value = '苏州·ICAN summer World动漫品牌夏游节'
$ws4.Cells.Item(34, 4).Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws4.Cells.Item(34, 5).Value = '2024.08.17 10:00-08.18 17:00'
$ws4.Cells.Item(34, 6).Value = 11331
$ws4.Cells.Item(34, 7).Value = 60
$ws4.Cells.Item(34, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85289'
$ws4.Cells.Item(34, 9).Value = '//i0.hdslb.com/bfs/openplatform/202404/JavlW9fj1714459472747.jpeg'

# row 35: 苏州·第二届Redamancy动漫游戏嘉年华
$ws4.Cells.Item(35, 2).NumberFormat = "@"
$ws4.Cells.Item(35, 2).Value = '2024-08-17'
$ws4.Cells.Item(35, 2).ClearFormats()
$ws4.Cells.Item(35, 3).Value = '苏州·第二届Redamancy动漫游戏嘉年华'
$ws4.Cells.Item(35, 4).Value = '清禾路886号 尹山湖大剧院'
$ws4.Cells.Item(35, 5).Value = '2024.08.17 10:00-08.18 17:00'
$ws4.Cells.Item(35, 6).Value = 1247
$ws4.Cells.Item(35, 7).Value = 60
$ws4.Cells.Item(35, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83576'
$ws4.Cells.Item(35, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/UbwYg1jn1716516632746.jpeg'

# row 36: 苏州·赛马娘ONLY
$ws4.Cells.Item(36, 2).NumberFormat = "@"
$ws4.Cells.Item(36, 2).Value = '2024-08-24'
$ws4.Cells.Item(36, 2).ClearFormats()
$ws4.Cells.Item(36, 3).Value = '苏州·赛马娘ONLY'
$ws4.Cells.Item(36, 4).Value = '东苑路115-11号 苏苑饭店'
$ws4.Cells.Item(36, 5).Value = '2024.08.24 10:00-08.24 16:00'
$ws4.Cells.Item(36, 6).Value = 18
$ws4.Cells.Item(36, 7).Value = 60
$ws4.Cells.Item(36, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87619'
$ws4.Cells.Item(36, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/2seg6U5l1718215220516.jpeg'

# row 37: 苏州·Good jump ACG中秋嘉年华动漫国潮文化节
$ws4.Cells.Item(37, 2).NumberFormat = "@"
$ws4.Cells.Item(37, 2).Value = '2024-09-15'
$ws4.Cells.Item(37, 2).ClearFormats()
$ws4.Cells.Item(37, 3).Value = '苏州·Good jump ACG中秋嘉年华动漫国潮文化节'
$ws4.Cells.Item(37, 4).Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws4.Cells.Item(37, 5).Value = '2024.09.15 10:00-09.16 17:00'
$ws4.Cells.Item(37, 6).Value = 150
$ws4.Cells.Item(37, 7).Value = 60
$ws4.Cells.Item(37, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87120'
$ws4.Cells.Item(37, 9).Value = '//i0.hdslb.com/bfs/openplatform/202406/5Qp6CuZ11717828215727.jpeg'

# row 38: 苏州·I COME ACG动漫品牌博览会
$ws4.Cells.Item(38, 2).NumberFormat = "@"
$ws4.Cells.Item(38, 2).Value = '2024-10-01'
$ws4.Cells.Item(38, 2).ClearFormats()
$ws4.Cells.Item(38, 3).Value = '苏州·I COME ACG动漫品牌博览会'
$ws4.Cells.Item(38, 4).Value = '金山南路288号木渎影视城F2 苏州广电国际会展中心'
$ws4.Cells.Item(38, 5).Value = '2024.10.01 10:00-10.03 17:00'
$ws4.Cells.Item(38, 6).Value = 207
$ws4.Cells.Item(38, 7).Value = 60
$ws4.Cells.Item(38, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=87118'
$ws4.Cells.Item(38, 9).Value = '//i2.hdslb.com/bfs/openplatform/202406/SyK3KnNb1717829071115.jpeg'

# row 39: 苏州·第十三届理想乡动漫展-同人创作者大会
$ws4.Cells.Item(39, 2).NumberFormat = "@"
$ws4.Cells.Item(39, 2).Value = '2024-10-01'
$ws4.Cells.Item(39, 2).ClearFormats()
$ws4.Cells.Item(39, 3).Value = '苏州·第十三届理想乡动漫展-同人创作者大会'
$ws4.Cells.Item(39, 4).Value = '花桥经济开发区绿地大道1598号 花桥国际博览中心'
$ws4.Cells.Item(39, 5).Value = '2024.10.01 10:00-10.03 17:00'
$ws4.Cells.Item(39, 6).Value = 3838
$ws4.Cells.Item(39, 7).Value = 39
$ws4.Cells.Item(39, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83821'
$ws4.Cells.Item(39, 9).Value = '//i0.hdslb.com/bfs/openplatform/202404/OMtuTTFY1711958198579.jpeg'

# row 40: 苏州·明日方舟ONLY#2024~佑桑柔
$ws4.Cells.Item(40, 2).NumberFormat = "@"
$ws4.Cells.Item(40, 2).Value = '2024-10-02'
$ws4.Cells.Item(40, 2).ClearFormats()
$ws4.Cells.Item(40, 3).Value = '苏州·明日方舟ONLY#2024~佑桑柔'
$ws4.Cells.Item(40, 4).Value = '城际路21号 苏州汇融广场假日酒店'
$ws4.Cells.Item(40, 5).Value = '2024.10.02 10:00-10.02 17:00'
$ws4.Cells.Item(40, 6).Value = 270
$ws4.Cells.Item(40, 7).Value = 75
$ws4.Cells.Item(40, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=84046'
$ws4.Cells.Item(40, 9).Value = '//i2.hdslb.com/bfs/openplatform/202405/0VhIRprD1716344515303.jpeg'

# row 41: 苏州·第三届华盟国漫次元嘉年华
$ws4.Cells.Item(41, 2).NumberFormat = "@"
$ws4.Cells.Item(41, 2).Value = '2024-10-26'
$ws4.Cells.Item(41, 2).ClearFormats()
$ws4.Cells.Item(41, 3).Value = '苏州·第三届华盟国漫次元嘉年华'
$ws4.Cells.Item(41, 4).Value = '清禾路886号 苏州聚橙尹山湖大剧院'
$ws4.Cells.Item(41, 5).Value = '2024.10.26 10:00-10.27 17:00'
$ws4.Cells.Item(41, 6).Value = 74
$ws4.Cells.Item(41, 7).Value = 58
$ws4.Cells.Item(41, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=85767'
$ws4.Cells.Item(41, 9).Value = '//i1.hdslb.com/bfs/openplatform/202405/CqSYBZhQ1715846719965.jpeg'

